# Auto-generated Excel COM-interop script
# Applies numeric cell value updates to the Sheets tables per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 462.5
$ws.Range("J2").Value = 425
$ws.Range("L2").Value = 425
$ws.Range("N2").Value = -651
$ws.Range("H15").Value = 405.03775
$ws.Range("I15").Value = 405.03775
$ws.Range("K15").Value = 1215.11325
$ws.Range("M15").Value = -1046.11325
$ws.Range("H18").Value = 1203.4642
$ws.Range("I18").Value = 550.1739
$ws.Range("J18").Value = 4208.6
$ws.Range("K18").Value = 550.1739
$ws.Range("L18").Value = 4208.6
$ws.Range("M18").Value = -266.1739
$ws.Range("N18").Value = -4776.6
$ws.Range("H28").Value = 40329.77
$ws.Range("I28").Value = 53216.95
$ws.Range("J28").Value = 5350.2856
$ws.Range("K28").Value = 53216.95
$ws.Range("L28").Value = 5350.2856
$ws.Range("M28").Value = -52731.95
$ws.Range("N28").Value = -6320.2856
$ws.Range("H33").Value = 16717405
$ws.Range("J33").Value = 62501052
$ws.Range("L33").Value = 62501052
$ws.Range("N33").Value = -62501510
$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996
$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216
$ws.Range("H70").Value = 1459269
$ws.Range("I70").Value = 3402694.2
$ws.Range("K70").Value = 10208082.6
$ws.Range("M70").Value = -10207812.6
$ws.Range("H73").Value = 1459269
$ws.Range("I73").Value = 3402694.2
$ws.Range("K73").Value = 10208082.6
$ws.Range("M73").Value = -10207146.6
$ws.Range("H98").Value = 2315.353
$ws.Range("I98").Value = 2147.5625
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2147.5625
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -649.5625
$ws.Range("N98").Value = -7996
$ws.Range("H101").Value = 1296.1818
$ws.Range("I101").Value = 280.5
$ws.Range("J101").Value = 2142.5833
$ws.Range("K101").Value = 841.5
$ws.Range("L101").Value = 6427.749899999999
$ws.Range("M101").Value = 780.5
$ws.Range("N101").Value = -9671.749899999999
$ws.Range("H106").Value = 102459.2
$ws.Range("I106").Value = 2735.75
$ws.Range("K106").Value = 2735.75
$ws.Range("M106").Value = -2104.75
$ws.Range("H107").Value = 342.72
$ws.Range("I107").Value = 356.7619
$ws.Range("J107").Value = 269
$ws.Range("K107").Value = 356.7619
$ws.Range("L107").Value = 269
$ws.Range("M107").Value = 1563.2381
$ws.Range("N107").Value = -4109
$ws.Range("H111").Value = 14899.305
$ws.Range("I111").Value = 527.2778
$ws.Range("J111").Value = 66638.60000000001
$ws.Range("K111").Value = 1581.8334
$ws.Range("L111").Value = 199915.8
$ws.Range("M111").Value = 1485.1666
$ws.Range("N111").Value = -206049.8
$ws.Range("H116").Value = 53142852
$ws.Range("J116").Value = 83338584
$ws.Range("L116").Value = 83338584
$ws.Range("N116").Value = -83345468
$ws.Range("H122").Value = 2315.353
$ws.Range("I122").Value = 2147.5625
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6442.6875
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3992.6875
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 4174.426
$ws.Range("I132").Value = 3672.3044
$ws.Range("K132").Value = 11016.9132
$ws.Range("M132").Value = -8486.913199999999
$ws.Range("H135").Value = 874.04
$ws.Range("I135").Value = 355.35294
$ws.Range("K135").Value = 3198.17646
$ws.Range("M135").Value = -663.1764599999997
$ws.Range("H137").Value = 1930.7106
$ws.Range("I137").Value = 1351.4073
$ws.Range("J137").Value = 3352.6365
$ws.Range("K137").Value = 4054.2219
$ws.Range("L137").Value = 10057.9095
$ws.Range("M137").Value = -1504.2219
$ws.Range("N137").Value = -15157.9095
$ws.Range("H138").Value = 2580.15
$ws.Range("I138").Value = 1958.875
$ws.Range("K138").Value = 5876.625
$ws.Range("M138").Value = -736.625
$ws.Range("H141").Value = 1529.6666
$ws.Range("I141").Value = 1181.2858
$ws.Range("J141").Value = 3968.3333
$ws.Range("K141").Value = 3543.8574
$ws.Range("L141").Value = 11904.9999
$ws.Range("M141").Value = 1636.1426
$ws.Range("N141").Value = -22264.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1204.8572
$ws.Range("I2").Value = 979.125
$ws.Range("K2").Value = 979.125
$ws.Range("M2").Value = -866.125
$ws.Range("H5").Value = 7731.5
$ws.Range("I5").Value = 784
$ws.Range("K5").Value = 784
$ws.Range("M5").Value = -672
$ws.Range("H32").Value = 7795.5845
$ws.Range("I32").Value = 2665.7144
$ws.Range("J32").Value = 21475.238
$ws.Range("K32").Value = 2665.7144
$ws.Range("L32").Value = 21475.238
$ws.Range("M32").Value = -2378.7144
$ws.Range("N32").Value = -22049.238
$ws.Range("H45").Value = 989.2
$ws.Range("I45").Value = 989.2
$ws.Range("K45").Value = 989.2
$ws.Range("M45").Value = -612.2
$ws.Range("H63").Value = 85722710
$ws.Range("I63").Value = 166676500
$ws.Range("J63").Value = 25007376
$ws.Range("K63").Value = 166676500
$ws.Range("L63").Value = 25007376
$ws.Range("M63").Value = -166675814
$ws.Range("N63").Value = -25008748
$ws.Range("H66").Value = 85722710
$ws.Range("I66").Value = 166676500
$ws.Range("J66").Value = 25007376
$ws.Range("K66").Value = 833382500
$ws.Range("L66").Value = 125036880
$ws.Range("M66").Value = -833379068
$ws.Range("N66").Value = -125043744
$ws.Range("H102").Value = 102884.125
$ws.Range("I102").Value = 90474.25
$ws.Range("K102").Value = 90474.25
$ws.Range("M102").Value = -88852.25
$ws.Range("H116").Value = 1204.8572
$ws.Range("I116").Value = 979.125
$ws.Range("K116").Value = 979.125
$ws.Range("M116").Value = 1314.875
$ws.Range("H122").Value = 4985.5713
$ws.Range("I122").Value = 4949.5
$ws.Range("K122").Value = 14848.5
$ws.Range("M122").Value = -12398.5
$ws.Range("H132").Value = 5037.143
$ws.Range("I132").Value = 3300
$ws.Range("J132").Value = 7353.3335
$ws.Range("K132").Value = 9900
$ws.Range("L132").Value = 22060.0005
$ws.Range("M132").Value = -7370
$ws.Range("N132").Value = -27120.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1204.8572
$ws.Range("I3").Value = 979.125
$ws.Range("K3").Value = 979.125
$ws.Range("M3").Value = -865.125
$ws.Range("H4").Value = 7731.5
$ws.Range("I4").Value = 784
$ws.Range("K4").Value = 784
$ws.Range("M4").Value = -669
$ws.Range("H20").Value = 10589
$ws.Range("I20").Value = 11408.538
$ws.Range("K20").Value = 11408.538
$ws.Range("M20").Value = -11161.538
$ws.Range("H86").Value = 2877.6428
$ws.Range("J86").Value = 3699.25
$ws.Range("L86").Value = 3699.25
$ws.Range("N86").Value = -5945.25
$ws.Range("H89").Value = 2877.6428
$ws.Range("J89").Value = 3699.25
$ws.Range("L89").Value = 18496.25
$ws.Range("N89").Value = -29728.25
$ws.Range("H94").Value = 19234800
$ws.Range("I94").Value = 41670724
$ws.Range("K94").Value = 41670724
$ws.Range("M94").Value = -41670273
$ws.Range("H105").Value = 2082.8
$ws.Range("I105").Value = 1750.5625
$ws.Range("J105").Value = 2673.4443
$ws.Range("K105").Value = 1750.5625
$ws.Range("L105").Value = 2673.4443
$ws.Range("M105").Value = -3.5625
$ws.Range("N105").Value = -6167.4443
$ws.Range("H107").Value = 33334122
$ws.Range("I107").Value = 661.1667
$ws.Range("J107").Value = 166667970
$ws.Range("K107").Value = 661.1667
$ws.Range("L107").Value = 166667970
$ws.Range("M107").Value = 1258.8333
$ws.Range("N107").Value = -166671810
$ws.Range("H134").Value = 3272.818
$ws.Range("I134").Value = 3182.7222
$ws.Range("J134").Value = 3678.25
$ws.Range("K134").Value = 9548.1666
$ws.Range("L134").Value = 11034.75
$ws.Range("M134").Value = -7013.1666
$ws.Range("N134").Value = -16104.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4967.591
$ws.Range("I31").Value = 15371.556
$ws.Range("K31").Value = 15371.556
$ws.Range("M31").Value = -15076.556
$ws.Range("H34").Value = 4967.591
$ws.Range("I34").Value = 15371.556
$ws.Range("K34").Value = 15371.556
$ws.Range("M34").Value = -15169.556
$ws.Range("H44").Value = 3001
$ws.Range("I44").Value = 3001
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3001
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -2559
$ws.Range("N44").ClearContents()
$ws.Range("H62").Value = 75957.14
$ws.Range("I62").Value = 4175
$ws.Range("K62").Value = 4175
$ws.Range("M62").Value = -3551
$ws.Range("H65").Value = 75957.14
$ws.Range("I65").Value = 4175
$ws.Range("K65").Value = 20875
$ws.Range("M65").Value = -17755
$ws.Range("H86").Value = 172099.8
$ws.Range("I86").Value = 416749.5
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 416749.5
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -415626.5
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 172099.8
$ws.Range("I89").Value = 416749.5
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 2083747.5
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -2078131.5
$ws.Range("N89").Value = -56232
$ws.Range("H94").Value = 1574.3636
$ws.Range("I94").Value = 1349.5
$ws.Range("J94").Value = 1702.8572
$ws.Range("K94").Value = 1349.5
$ws.Range("L94").Value = 1702.8572
$ws.Range("M94").Value = -898.5
$ws.Range("N94").Value = -2604.8572
$ws.Range("H99").Value = 2931.25
$ws.Range("I99").Value = 2590
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2590
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1092
$ws.Range("N99").Value = -6496
$ws.Range("H107").Value = 1201.5
$ws.Range("I107").Value = 678.25
$ws.Range("K107").Value = 678.25
$ws.Range("M107").Value = 1241.75
$ws.Range("H122").Value = 2249.3684
$ws.Range("I122").Value = 2192.375
$ws.Range("J122").Value = 2290.818
$ws.Range("K122").Value = 6577.125
$ws.Range("L122").Value = 6872.454000000001
$ws.Range("M122").Value = -4127.125
$ws.Range("N122").Value = -11772.454
$ws.Range("H126").Value = 2931.25
$ws.Range("I126").Value = 2590
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7770
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5300
$ws.Range("N126").Value = -15440
$ws.Range("H134").Value = 3123.16
$ws.Range("I134").Value = 3003.3333
$ws.Range("K134").Value = 9009.999899999999
$ws.Range("M134").Value = -6474.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 873.5
$ws.Range("I2").Value = 812.7857
$ws.Range("J2").Value = 944.3333
$ws.Range("K2").Value = 4876.7142
$ws.Range("L2").Value = 5665.9998
$ws.Range("M2").Value = -4763.7142
$ws.Range("N2").Value = -5891.9998
$ws.Range("H34").Value = 1733.2307
$ws.Range("J34").Value = 2618.5
$ws.Range("L34").Value = 7855.5
$ws.Range("N34").Value = -8023.5
$ws.Range("H47").Value = 6801
$ws.Range("I47").Value = 403
$ws.Range("K47").Value = 1209
$ws.Range("M47").Value = -778
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H51").Value = 1018
$ws.Range("I51").Value = 852
$ws.Range("J51").Value = 1350
$ws.Range("K51").Value = 2556
$ws.Range("L51").Value = 4050
$ws.Range("M51").Value = -2096
$ws.Range("N51").Value = -4970
$ws.Range("H57").Value = 500
$ws.Range("I57").Value = 500
$ws.Range("K57").Value = 1500
$ws.Range("M57").Value = -941
$ws.Range("H58").Value = 3550
$ws.Range("I58").Value = 825
$ws.Range("K58").Value = 2475
$ws.Range("M58").Value = -2347
$ws.Range("H97").Value = 659.1
$ws.Range("I97").Value = 567.6923
$ws.Range("J97").Value = 828.8570999999999
$ws.Range("K97").Value = 1703.0769
$ws.Range("L97").Value = 2486.5713
$ws.Range("M97").Value = -1207.0769
$ws.Range("N97").Value = -3478.5713
$ws.Range("H105").Value = 15265.733
$ws.Range("H114").Value = 4206.615
$ws.Range("I114").Value = 1431.1666
$ws.Range("J114").Value = 6585.5713
$ws.Range("K114").Value = 4293.4998
$ws.Range("L114").Value = 19756.7139
$ws.Range("M114").Value = -1039.4998
$ws.Range("N114").Value = -26264.7139
$ws.Range("H117").Value = 391.16666
$ws.Range("J117").Value = 1200
$ws.Range("L117").Value = 3600
$ws.Range("N117").Value = -10484
$ws.Range("H121").Value = 361578.34
$ws.Range("J121").Value = 722.5833
$ws.Range("L121").Value = 2167.7499
$ws.Range("N121").Value = -4787.7499
$ws.Range("H132").Value = 2692.7
$ws.Range("I132").Value = 1050
$ws.Range("J132").Value = 2766.2537
$ws.Range("K132").Value = 9450
$ws.Range("L132").Value = 24896.2833
$ws.Range("M132").Value = -6920
$ws.Range("N132").Value = -29956.2833

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H122").Value = 2292.2632
$ws.Range("I122").Value = 1840.8667
$ws.Range("J122").Value = 3985
$ws.Range("K122").Value = 5522.6001
$ws.Range("L122").Value = 11955
$ws.Range("M122").Value = -3072.6001
$ws.Range("N122").Value = -16855
$ws.Range("H132").Value = 7492.2666
$ws.Range("I132").Value = 6948.75
$ws.Range("J132").Value = 8113.4287
$ws.Range("K132").Value = 20846.25
$ws.Range("L132").Value = 24340.2861
$ws.Range("M132").Value = -18316.25
$ws.Range("N132").Value = -29400.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9619.556
$ws.Range("I61").Value = 11471
$ws.Range("J61").Value = 3139.5
$ws.Range("K61").Value = 11471
$ws.Range("L61").Value = 3139.5
$ws.Range("M61").Value = -11269
$ws.Range("N61").Value = -3543.5
$ws.Range("H68").Value = 10228.777
$ws.Range("I68").Value = 4921
$ws.Range("J68").Value = 20844.334
$ws.Range("K68").Value = 4921
$ws.Range("L68").Value = 20844.334
$ws.Range("M68").Value = -4172
$ws.Range("N68").Value = -22342.334
$ws.Range("H71").Value = 10228.777
$ws.Range("I71").Value = 4921
$ws.Range("J71").Value = 20844.334
$ws.Range("K71").Value = 24605
$ws.Range("L71").Value = 104221.67
$ws.Range("M71").Value = -20861
$ws.Range("N71").Value = -111709.67
$ws.Range("H93").Value = 25644832
$ws.Range("J93").Value = 4199.5
$ws.Range("L93").Value = 4199.5
$ws.Range("N93").Value = -6695.5
$ws.Range("H113").Value = 9619.556
$ws.Range("I113").Value = 11471
$ws.Range("J113").Value = 3139.5
$ws.Range("K113").Value = 11471
$ws.Range("L113").Value = 3139.5
$ws.Range("M113").Value = -9301
$ws.Range("N113").Value = -7479.5
$ws.Range("H122").Value = 4601.2573
$ws.Range("I122").Value = 2999.5652
$ws.Range("K122").Value = 8998.695599999999
$ws.Range("M122").Value = -6548.695599999999
$ws.Range("H132").Value = 18522476
$ws.Range("I132").Value = 28574010
$ws.Range("J132").Value = 6490.684
$ws.Range("K132").Value = 85722030
$ws.Range("L132").Value = 19472.052
$ws.Range("M132").Value = -85719500
$ws.Range("N132").Value = -24532.052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22224358
$ws.Range("J81").Value = 200000000
$ws.Range("L81").Value = 400000000
$ws.Range("N81").Value = -400002122
$ws.Range("H84").Value = 22224358
$ws.Range("J84").Value = 200000000
$ws.Range("L84").Value = 2000000000
$ws.Range("N84").Value = -2000010608
$ws.Range("H107").Value = 698.6
$ws.Range("I107").Value = 533.4706
$ws.Range("K107").Value = 1600.4118
$ws.Range("M107").Value = 319.5882000000001
$ws.Range("H113").Value = 6741.769
$ws.Range("I113").Value = 9374.5
$ws.Range("J113").Value = 4485.143
$ws.Range("K113").Value = 28123.5
$ws.Range("L113").Value = 13455.429
$ws.Range("M113").Value = -25953.5
$ws.Range("N113").Value = -17795.429
$ws.Range("H126").Value = 1813.0714
$ws.Range("I126").Value = 1487.1111
$ws.Range("J126").Value = 2399.8
$ws.Range("K126").Value = 4461.3333
$ws.Range("L126").Value = 7199.400000000001
$ws.Range("M126").Value = -1991.3333
$ws.Range("N126").Value = -12139.4
$ws.Range("H132").Value = 404444.66
$ws.Range("I132").Value = 516485.78
$ws.Range("J132").Value = 7208
$ws.Range("K132").Value = 1549457.34
$ws.Range("L132").Value = 21624
$ws.Range("M132").Value = -1546927.34
$ws.Range("N132").Value = -26684
$ws.Range("H136").Value = 1966.0698
$ws.Range("I136").Value = 1674.7894
$ws.Range("K136").Value = 5024.3682
$ws.Range("M136").Value = -2474.3682
